# Update the campaign dates heading text across the document (it occurs
# multiple times, once per language/section block). Word's Find/Replace
# with Replace:=2 (wdReplaceAll) will update every matching occurrence.

$d = $word.ActiveDocument

$find = "Kampagnendaten Sternbild Stier 2022: 16.-25. Januar"
$replace = "Kampagnendaten 2022 für das Sternbild Sternbild Stier: 16.-25. Januar"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
